# Auto-generated Excel COM-interop script
# Applies cell-level numeric corrections to the Leve profit tables
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 342
$ws.Range("I55").Value = 342
$ws.Range("K55").Value = 342
$ws.Range("M55").Value = -128
$ws.Range("H98").Value = 1876.7142
$ws.Range("I98").Value = 1876.7142
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1876.7142
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -378.7141999999999
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1876.7142
$ws.Range("I122").Value = 1876.7142
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5630.142599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3180.142599999999
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 2131.6667
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 2697.5
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 24277.5
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -29197.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2052.2104
$ws.Range("I2").Value = 1547.909
$ws.Range("J2").Value = 2745.625
$ws.Range("K2").Value = 1547.909
$ws.Range("L2").Value = 2745.625
$ws.Range("M2").Value = -1434.909
$ws.Range("N2").Value = -2971.625
$ws.Range("H61").Value = 2700
$ws.Range("I61").Value = 2700
$ws.Range("K61").Value = 2700
$ws.Range("M61").Value = -2488
$ws.Range("H110").Value = 1444.3077
$ws.Range("I110").Value = 1539.6666
$ws.Range("K110").Value = 1539.6666
$ws.Range("M110").Value = 505.3334
$ws.Range("H116").Value = 2052.2104
$ws.Range("I116").Value = 1547.909
$ws.Range("J116").Value = 2745.625
$ws.Range("K116").Value = 1547.909
$ws.Range("L116").Value = 2745.625
$ws.Range("M116").Value = 746.0909999999999
$ws.Range("N116").Value = -7333.625
$ws.Range("H122").Value = 3026.4443
$ws.Range("I122").Value = 2539.8333
$ws.Range("K122").Value = 7619.499899999999
$ws.Range("M122").Value = -5169.499899999999
$ws.Range("H136").Value = 2700
$ws.Range("I136").Value = 2700
$ws.Range("K136").Value = 8100
$ws.Range("M136").Value = -5550

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2052.2104
$ws.Range("I3").Value = 1547.909
$ws.Range("J3").Value = 2745.625
$ws.Range("K3").Value = 1547.909
$ws.Range("L3").Value = 2745.625
$ws.Range("M3").Value = -1433.909
$ws.Range("N3").Value = -2973.625
$ws.Range("H64").Value = 858.2
$ws.Range("I64").Value = 796.3333
$ws.Range("J64").Value = 951
$ws.Range("K64").Value = 796.3333
$ws.Range("L64").Value = 951
$ws.Range("M64").Value = -571.3333
$ws.Range("N64").Value = -1401
$ws.Range("H67").Value = 858.2
$ws.Range("I67").Value = 796.3333
$ws.Range("J67").Value = 951
$ws.Range("K67").Value = 796.3333
$ws.Range("L67").Value = 951
$ws.Range("M67").Value = -16.33330000000001
$ws.Range("N67").Value = -2511
$ws.Range("H94").Value = 1537.8334
$ws.Range("I94").Value = 1404.9546
$ws.Range("K94").Value = 1404.9546
$ws.Range("M94").Value = -953.9546
$ws.Range("H102").Value = 64338.75
$ws.Range("I102").Value = 64338.75
$ws.Range("K102").Value = 64338.75
$ws.Range("M102").Value = -61093.75
$ws.Range("H134").Value = 1388.1428
$ws.Range("I134").Value = 1202.8334
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 3608.5002
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -1073.5002
$ws.Range("N134").Value = -12570

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5258.125
$ws.Range("I132").Value = 5258.125
$ws.Range("K132").Value = 15774.375
$ws.Range("M132").Value = -13244.375
$ws.Range("H134").Value = 2266
$ws.Range("I134").Value = 2184.4443
$ws.Range("K134").Value = 6553.3329
$ws.Range("M134").Value = -4018.3329

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 16000
$ws.Range("J94").Value = 19000
$ws.Range("L94").Value = 57000
$ws.Range("N94").Value = -58352
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H113").Value = 201191
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 251113.75
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 753341.25
$ws.Range("M113").Value = -2330
$ws.Range("N113").Value = -757681.25
$ws.Range("H121").Value = 16673
$ws.Range("I121").Value = 22511.6
$ws.Range("J121").Value = 9374.75
$ws.Range("K121").Value = 67534.79999999999
$ws.Range("L121").Value = 28124.25
$ws.Range("M121").Value = -66224.79999999999
$ws.Range("N121").Value = -30744.25
$ws.Range("H134").Value = 1229.5
$ws.Range("I134").Value = 1106.1666
$ws.Range("K134").Value = 3318.4998
$ws.Range("M134").Value = 1751.5002

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3999
$ws.Range("J113").Value = 3999
$ws.Range("L113").Value = 3999
$ws.Range("N113").Value = -8339
$ws.Range("H122").Value = 3559.25
$ws.Range("I122").Value = 3559.25
$ws.Range("K122").Value = 10677.75
$ws.Range("M122").Value = -8227.75
$ws.Range("H132").Value = 4918.15
$ws.Range("I132").Value = 5349.2354
$ws.Range("J132").Value = 2475.3333
$ws.Range("K132").Value = 16047.7062
$ws.Range("L132").Value = 7425.999899999999
$ws.Range("M132").Value = -13517.7062
$ws.Range("N132").Value = -12485.9999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6723.909
$ws.Range("H46").Value = 2744.9092
$ws.Range("J46").Value = 2631.2856
$ws.Range("L46").Value = 2631.2856
$ws.Range("N46").Value = -3007.2856
$ws.Range("H100").Value = 4043.4375
$ws.Range("I100").Value = 4330.6924
$ws.Range("J100").Value = 2798.6667
$ws.Range("K100").Value = 4330.6924
$ws.Range("L100").Value = 2798.6667
$ws.Range("M100").Value = -3789.6924
$ws.Range("N100").Value = -3880.6667
$ws.Range("H122").Value = 6870.6875
$ws.Range("J122").Value = 8099
$ws.Range("L122").Value = 24297
$ws.Range("N122").Value = -29197
$ws.Range("H126").Value = 6723.909
$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 39974.5
$ws.Range("J15").Value = 39974.5
$ws.Range("L15").Value = 39974.5
$ws.Range("N15").Value = -40550.5
$ws.Range("H101").Value = 11863.667
$ws.Range("J101").Value = 11863.667
$ws.Range("L101").Value = 11863.667
$ws.Range("N101").Value = -18353.667
$ws.Range("H122").Value = 6158.3335
$ws.Range("I122").Value = 5284.6665
$ws.Range("K122").Value = 15853.9995
$ws.Range("M122").Value = -13403.9995
$ws.Range("H126").Value = 1661.375
$ws.Range("I126").Value = 1248.8334
$ws.Range("K126").Value = 3746.5002
$ws.Range("M126").Value = -1276.5002
$ws.Range("H136").Value = 4273.5
$ws.Range("I136").Value = 4273.5
$ws.Range("K136").Value = 12820.5
$ws.Range("M136").Value = -10270.5
